$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.274.11'
$ws.Cells.Item(2, 5).Value = '  +3.01%  '

$ws.Cells.Item(3, 4).Value = '1.902.06'
$ws.Cells.Item(3, 5).Value = '  +1.47%  '

$ws.Cells.Item(4, 4).Value = '1.004'
$ws.Cells.Item(4, 5).Value = '  -1.63%  '

$ws.Cells.Item(5, 4).Value = '315.35'
$ws.Cells.Item(5, 5).Value = '  -0.50%  '

$ws.Cells.Item(6, 5).Value = '  -1.44%  '

$ws.Cells.Item(7, 4).Value = '0.5124'
$ws.Cells.Item(7, 5).Value = '  +0.25%  '

$ws.Cells.Item(8, 4).Value = '0.3942'
$ws.Cells.Item(8, 5).Value = '  -0.09%  '

$ws.Cells.Item(9, 4).Value = '0.08455'
$ws.Cells.Item(9, 5).Value = '  +0.15%  '

$ws.Cells.Item(10, 4).Value = '42.51'
$ws.Cells.Item(10, 5).Value = '  +1.27%  '

$ws.Cells.Item(11, 5).Value = '  +1.13%  '

$ws.Cells.Item(12, 2).Value = 'Polkadot'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(12, 4).Value = '6.254'
$ws.Cells.Item(12, 5).Value = '  +0.22%  '

$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).Value = '1.901.38'
$ws.Cells.Item(13, 5).Value = '  +1.18%  '

$ws.Cells.Item(14, 4).Value = '20.64'
$ws.Cells.Item(14, 5).Value = '  +0.99%  '

$ws.Cells.Item(15, 4).Value = '7.365'
$ws.Cells.Item(15, 5).Value = '  +1.84%  '

$ws.Cells.Item(16, 4).Value = '1.004'
$ws.Cells.Item(16, 5).Value = '  -1.73%  '

$ws.Cells.Item(17, 4).Value = '93.18'
$ws.Cells.Item(17, 5).Value = '  +2.46%  '

$ws.Cells.Item(18, 4).Value = '0.00001106'
$ws.Cells.Item(18, 5).Value = '  -0.19%  '

$ws.Cells.Item(19, 4).Value = '0.06725'
$ws.Cells.Item(19, 5).Value = '  -0.67%  '

$ws.Cells.Item(20, 4).Value = '17.92'
$ws.Cells.Item(20, 5).Value = '  +1.35%  '

$ws.Cells.Item(21, 5).Value = '  -1.42%  '

$ws.Cells.Item(22, 4).Value = '6.039'
$ws.Cells.Item(22, 5).Value = '  +1.77%  '

$ws.Cells.Item(23, 4).Value = '29.286.31'
$ws.Cells.Item(23, 5).Value = '  +2.88%  '

$ws.Cells.Item(24, 4).Value = '11.17'
$ws.Cells.Item(24, 5).Value = '  +0.18%  '

$ws.Cells.Item(25, 5).Value = '  -3.30%  '

$ws.Cells.Item(26, 4).Value = '2.119.97'
$ws.Cells.Item(26, 5).Value = '  +1.41%  '

$ws.Cells.Item(27, 4).Value = '160.30'
$ws.Cells.Item(27, 5).Value = '  -0.90%  '

$ws.Cells.Item(28, 4).Value = '20.95'
$ws.Cells.Item(28, 5).Value = '  +1.03%  '

$ws.Cells.Item(29, 4).Value = '2.440'
$ws.Cells.Item(29, 5).Value = '  +4.03%  '

$ws.Cells.Item(30, 4).Value = '126.99'
$ws.Cells.Item(30, 5).Value = '  -0.10%  '

$ws.Cells.Item(31, 4).Value = '1.060'
$ws.Cells.Item(31, 5).Value = '  +1.93%  '

$ws.Cells.Item(32, 4).Value = '0.1046'
$ws.Cells.Item(32, 5).Value = '  -0.83%  '

$ws.Cells.Item(33, 4).Value = '5.982'
$ws.Cells.Item(33, 5).Value = '  +4.08%  '

$ws.Cells.Item(34, 4).Value = '3.649'
$ws.Cells.Item(34, 5).Value = '  +0.15%  '

$ws.Cells.Item(35, 4).Value = '0.02477'
$ws.Cells.Item(35, 5).Value = '  +1.84%  '

$ws.Cells.Item(36, 4).Value = '0.06616'
$ws.Cells.Item(36, 5).Value = '  +2.34%  '

$ws.Cells.Item(37, 4).Value = '9.126'
$ws.Cells.Item(37, 5).Value = '  +3.78%  '

$ws.Cells.Item(38, 4).Value = '0.2196'
$ws.Cells.Item(38, 5).Value = '  +1.09%  '

$ws.Cells.Item(39, 4).Value = '1.234'
$ws.Cells.Item(39, 5).Value = '  +4.10%  '

$ws.Cells.Item(40, 4).Value = '5.116'
$ws.Cells.Item(40, 5).Value = '  +2.60%  '

$ws.Cells.Item(41, 4).Value = '0.6490'
$ws.Cells.Item(41, 5).Value = '  +1.78%  '

$ws.Cells.Item(42, 4).Value = '1.236'
$ws.Cells.Item(42, 5).Value = '  -2.35%  '

$ws.Cells.Item(43, 4).Value = '11.28'
$ws.Cells.Item(43, 5).Value = '  +0.74%  '

$ws.Cells.Item(44, 5).Value = '  -1.43%  '

$ws.Cells.Item(45, 4).Value = '0.6053'
$ws.Cells.Item(45, 5).Value = '  +0.27%  '

$ws.Cells.Item(46, 5).Value = '  +1.09%  '

$ws.Cells.Item(47, 4).Value = '3.680'
$ws.Cells.Item(47, 5).Value = '  -0.82%  '

$ws.Cells.Item(48, 4).Value = '2.052'
$ws.Cells.Item(48, 5).Value = '  +3.18%  '

$ws.Cells.Item(49, 4).Value = '1.231'
$ws.Cells.Item(49, 5).Value = '  +2.25%  '

$ws.Cells.Item(50, 4).Value = '123.10'
$ws.Cells.Item(50, 5).Value = '  +0.97%  '

$ws.Cells.Item(51, 4).Value = '1.181'
$ws.Cells.Item(51, 5).Value = '  -2.09%  '
